# Cuttings log well 15-6-10.xlsx — apply the commit's edits via Excel COM.
#
# Summary of the underlying edit (per the OOXML diff):
#  - The "Simple lithology" classification formula in column D was extended
#    so that "Marl" rock is now bucketed into "Shale" (previously only
#    Claystone/Siltstone/Shale/Clay were folded into "Shale"). This is the
#    substantive, data-affecting change and ripples through every row whose
#    lithology is Marl (their D value flips from "Marl" to "Shale").
#  - A duplicate shared-string entry ("coal", lower-case) was removed in
#    favour of the existing "Coal" entry — row 135's lithology cell now
#    points at "Coal" instead of the stray lower-case duplicate.
#  - The active selection/scroll position left behind by the editing user
#    moved from D2 to C136 (with the view scrolled down near row 118).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the duplicated shared string: row 135's lithology was "coal"
#        (lower case); make it match the canonical "Coal" used elsewhere
#        (e.g. row 138). Once nothing references "coal" any more, the
#        exporter drops it from the shared-strings table automatically.
$ws.Range("C135").Value = "Coal"

# --- 2. Extend the lithology-bucketing formula to also fold "Marl" into
#        "Shale". The sheet uses three shared-formula groups in column D
#        (D3:D66, D67:D130, D131:D155) plus one standalone cell (D2);
#        assigning .Formula across each *whole* range at once (rather than
#        cell-by-cell) keeps Excel's shared-formula grouping intact instead
#        of exploding it into one literal formula per cell.
$ws.Range("D2").Formula = '=IF(OR(C2="Claystone",C2="Siltstone",C2="Shale",C2="Clay",C2="Marl"),"Shale",IF(C2="Limestone","Carbonate",C2))'
$ws.Range("D3:D66").Formula = '=IF(OR(C3="Claystone",C3="Siltstone",C3="Shale",C3="Clay",C3="Marl"),"Shale",IF(C3="Limestone","Carbonate",C3))'
$ws.Range("D67:D130").Formula = '=IF(OR(C67="Claystone",C67="Siltstone",C67="Shale",C67="Clay",C67="Marl"),"Shale",IF(C67="Limestone","Carbonate",C67))'
$ws.Range("D131:D155").Formula = '=IF(OR(C131="Claystone",C131="Siltstone",C131="Shale",C131="Clay",C131="Marl"),"Shale",IF(C131="Limestone","Carbonate",C131))'

# --- 3. Move the saved selection/scroll position to where the user left it:
#        cell C136, with the window scrolled so row 118 is near the top.
$win = $excel.ActiveWindow
$win.ScrollRow = 118
$win.ScrollColumn = 1
$ws.Range("C136").Select()
